$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recompute ligand/receptor-expressing cell counts (E, K: 1 -> 3) and the
# downstream total expression, specificity, and edge-weight columns that
# depend on them, per Dr Hou advice.

$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = 3
$row2[0,1] = 1
$row2[0,2] = 8.510923
$row2[0,3] = 25.532769
$row2[0,4] = 0.028708534116067
$row2[0,5] = 0.028708534116067
$row2[0,6] = 3
$row2[0,7] = 1
$row2[0,8] = 23.18520366666667
$row2[0,9] = 69.555611
$row2[0,10] = 0.4216200689608106
$row2[0,11] = 0.4216200689608105
$row2[0,12] = 197.3274831463177
$row2[0,13] = 1775.947348316859
$row2[0,14] = 0.01210409413377995
$row2[0,15] = 0.01210409413377995
$ws.Range("E2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = 3
$row3[0,1] = 1
$row3[0,2] = 8.510923
$row3[0,3] = 25.532769
$row3[0,4] = 0.028708534116067
$row3[0,5] = 0.028708534116067
$row3[0,6] = 3
$row3[0,7] = 1
$row3[0,8] = 11.56543033333333
$row3[0,9] = 34.696291
$row3[0,10] = 0.210315924104302
$row3[0,11] = 0.2103159241043019
$row3[0,12] = 98.43248702886433
$row3[0,13] = 885.8923832597791
$row3[0,14] = 0.006037861882300512
$row3[0,15] = 0.006037861882300511
$ws.Range("E3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = 3
$row4[0,1] = 1
$row4[0,2] = 8.510923
$row4[0,3] = 25.532769
$row4[0,4] = 0.028708534116067
$row4[0,5] = 0.028708534116067
$row4[0,6] = 3
$row4[0,7] = 1
$row4[0,8] = 2.096289333333333
$row4[0,9] = 6.288868
$row4[0,10] = 0.03812076296541245
$row4[0,11] = 0.03812076296541244
$row4[0,12] = 17.84135710172134
$row4[0,13] = 160.572213915492
$row4[0,14] = 0.001094391224123047
$row4[0,15] = 0.001094391224123047
$ws.Range("E4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = 3
$row5[0,1] = 1
$row5[0,2] = 8.510923
$row5[0,3] = 25.532769
$row5[0,4] = 0.028708534116067
$row5[0,5] = 0.028708534116067
$row5[0,6] = 3
$row5[0,7] = 1
$row5[0,8] = 18.14382633333333
$row5[0,9] = 54.431479
$row5[0,10] = 0.3299432439694752
$row5[0,11] = 0.3299432439694752
$row5[0,12] = 154.4207088483723
$row5[0,13] = 1389.786379635351
$row5[0,14] = 0.009472186875863495
$row5[0,15] = 0.009472186875863497
$ws.Range("E5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = 3
$row6[0,1] = 1
$row6[0,2] = 251.3975576666667
$row6[0,3] = 754.192673
$row6[0,4] = 0.847999137222769
$row6[0,5] = 0.8479991372227691
$row6[0,6] = 3
$row6[0,7] = 1
$row6[0,8] = 23.18520366666667
$row6[0,9] = 69.555611
$row6[0,10] = 0.4216200689608106
$row6[0,11] = 0.4216200689608105
$row6[0,12] = 5828.703575804245
$row6[0,13] = 52458.33218223821
$row6[0,14] = 0.3575334547145717
$row6[0,15] = 0.3575334547145717
$ws.Range("E6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = 3
$row7[0,1] = 1
$row7[0,2] = 251.3975576666667
$row7[0,3] = 754.192673
$row7[0,4] = 0.847999137222769
$row7[0,5] = 0.8479991372227691
$row7[0,6] = 3
$row7[0,7] = 1
$row7[0,8] = 11.56543033333333
$row7[0,9] = 34.696291
$row7[0,10] = 0.210315924104302
$row7[0,11] = 0.2103159241043019
$row7[0,12] = 2907.520939163983
$row7[0,13] = 26167.68845247585
$row7[0,14] = 0.1783477221846574
$row7[0,15] = 0.1783477221846574
$ws.Range("E7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,16
$row8[0,0] = 3
$row8[0,1] = 1
$row8[0,2] = 251.3975576666667
$row8[0,3] = 754.192673
$row8[0,4] = 0.847999137222769
$row8[0,5] = 0.8479991372227691
$row8[0,6] = 3
$row8[0,7] = 1
$row8[0,8] = 2.096289333333333
$row8[0,9] = 6.288868
$row8[0,10] = 0.03812076296541245
$row8[0,11] = 0.03812076296541244
$row8[0,12] = 527.002018562685
$row8[0,13] = 4743.018167064164
$row8[0,14] = 0.03232637410494345
$row8[0,15] = 0.03232637410494344
$ws.Range("E8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,16
$row9[0,0] = 3
$row9[0,1] = 1
$row9[0,2] = 251.3975576666667
$row9[0,3] = 754.192673
$row9[0,4] = 0.847999137222769
$row9[0,5] = 0.8479991372227691
$row9[0,6] = 3
$row9[0,7] = 1
$row9[0,8] = 18.14382633333333
$row9[0,9] = 54.431479
$row9[0,10] = 0.3299432439694752
$row9[0,11] = 0.3299432439694752
$row9[0,12] = 4561.313626928152
$row9[0,13] = 41051.82264235337
$row9[0,14] = 0.2797915862185965
$row9[0,15] = 0.2797915862185966
$ws.Range("E9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,16
$row10[0,0] = 3
$row10[0,1] = 1
$row10[0,2] = 19.800378
$row10[0,3] = 59.401134
$row10[0,4] = 0.06678944543664916
$row10[0,5] = 0.06678944543664918
$row10[0,6] = 3
$row10[0,7] = 1
$row10[0,8] = 23.18520366666667
$row10[0,9] = 69.555611
$row10[0,10] = 0.4216200689608106
$row10[0,11] = 0.4216200689608105
$row10[0,12] = 459.0757966069859
$row10[0,13] = 4131.682169462874
$row10[0,14] = 0.02815977059085431
$row10[0,15] = 0.02815977059085432
$ws.Range("E10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,16
$row11[0,0] = 3
$row11[0,1] = 1
$row11[0,2] = 19.800378
$row11[0,3] = 59.401134
$row11[0,4] = 0.06678944543664916
$row11[0,5] = 0.06678944543664918
$row11[0,6] = 3
$row11[0,7] = 1
$row11[0,8] = 11.56543033333333
$row11[0,9] = 34.696291
$row11[0,10] = 0.210315924104302
$row11[0,11] = 0.2103159241043019
$row11[0,12] = 228.999892332666
$row11[0,13] = 2060.999030993994
$row11[0,14] = 0.01404688393742272
$row11[0,15] = 0.01404688393742272
$ws.Range("E11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,16
$row12[0,0] = 3
$row12[0,1] = 1
$row12[0,2] = 19.800378
$row12[0,3] = 59.401134
$row12[0,4] = 0.06678944543664916
$row12[0,5] = 0.06678944543664918
$row12[0,6] = 3
$row12[0,7] = 1
$row12[0,8] = 2.096289333333333
$row12[0,9] = 6.288868
$row12[0,10] = 0.03812076296541245
$row12[0,11] = 0.03812076296541244
$row12[0,12] = 41.507321197368
$row12[0,13] = 373.565890776312
$row12[0,14] = 0.002546064618081851
$row12[0,15] = 0.002546064618081851
$ws.Range("E12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,16
$row13[0,0] = 3
$row13[0,1] = 1
$row13[0,2] = 19.800378
$row13[0,3] = 59.401134
$row13[0,4] = 0.06678944543664916
$row13[0,5] = 0.06678944543664918
$row13[0,6] = 3
$row13[0,7] = 1
$row13[0,8] = 18.14382633333333
$row13[0,9] = 54.431479
$row13[0,10] = 0.3299432439694752
$row13[0,11] = 0.3299432439694752
$row13[0,12] = 359.254619766354
$row13[0,13] = 3233.291577897186
$row13[0,14] = 0.02203672629029028
$row13[0,15] = 0.02203672629029029
$ws.Range("E13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,16
$row14[0,0] = 3
$row14[0,1] = 1
$row14[0,2] = 16.750827
$row14[0,3] = 50.252481
$row14[0,4] = 0.05650288322451468
$row14[0,5] = 0.0565028832245147
$row14[0,6] = 3
$row14[0,7] = 1
$row14[0,8] = 23.18520366666667
$row14[0,9] = 69.555611
$row14[0,10] = 0.4216200689608106
$row14[0,11] = 0.4216200689608105
$row14[0,12] = 388.371335580099
$row14[0,13] = 3495.342020220891
$row14[0,14] = 0.02382274952160451
$row14[0,15] = 0.02382274952160451
$ws.Range("E14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,16
$row15[0,0] = 3
$row15[0,1] = 1
$row15[0,2] = 16.750827
$row15[0,3] = 50.252481
$row15[0,4] = 0.05650288322451468
$row15[0,5] = 0.0565028832245147
$row15[0,6] = 3
$row15[0,7] = 1
$row15[0,8] = 11.56543033333333
$row15[0,9] = 34.696291
$row15[0,10] = 0.210315924104302
$row15[0,11] = 0.2103159241043019
$row15[0,12] = 193.730522694219
$row15[0,13] = 1743.574704247971
$row15[0,14] = 0.01188345609992127
$row15[0,15] = 0.01188345609992127
$ws.Range("E15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,16
$row16[0,0] = 3
$row16[0,1] = 1
$row16[0,2] = 16.750827
$row16[0,3] = 50.252481
$row16[0,4] = 0.05650288322451468
$row16[0,5] = 0.0565028832245147
$row16[0,6] = 3
$row16[0,7] = 1
$row16[0,8] = 2.096289333333333
$row16[0,9] = 6.288868
$row16[0,10] = 0.03812076296541245
$row16[0,11] = 0.03812076296541244
$row16[0,12] = 35.114579964612
$row16[0,13] = 316.031219681508
$row16[0,14] = 0.002153933018264104
$row16[0,15] = 0.002153933018264104
$ws.Range("E16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,16
$row17[0,0] = 3
$row17[0,1] = 1
$row17[0,2] = 16.750827
$row17[0,3] = 50.252481
$row17[0,4] = 0.05650288322451468
$row17[0,5] = 0.0565028832245147
$row17[0,6] = 3
$row17[0,7] = 1
$row17[0,8] = 18.14382633333333
$row17[0,9] = 54.431479
$row17[0,10] = 0.3299432439694752
$row17[0,11] = 0.3299432439694752
$row17[0,12] = 303.924096027711
$row17[0,13] = 2735.316864249399
$row17[0,14] = 0.01864274458472481
$row17[0,15] = 0.01864274458472482
$ws.Range("E17:T17").Value = $row17
